$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '87.677.62'
$ws.Range("E2").Value = '  +6.64%  '
$ws.Range("D3").Value = '3.303.50'
$ws.Range("E3").Value = '  +3.18%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '633.38'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.400'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +35.96%  '
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.647'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.26%  '
$ws.Range("B9").Value = 'USDC'
$ws.Range("C9").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '3.298.12'
$ws.Range("E10").Value = '  +3.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.594'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000268'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.177'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.76%  '
$ws.Range("D15").Value = '3.909.27'
$ws.Range("E16").Value = '  -0.37%  '
$ws.Range("D17").Value = '87.361.86'
$ws.Range("E17").Value = '  +6.63%  '
$ws.Range("D18").Value = '3.294.20'
$ws.Range("E18").Value = '  +3.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.33'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.01'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -8.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '440.90'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.44'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.27'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.40'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.51%  '
$ws.Range("E26").Value = '  -2.10%  '
$ws.Range("D27").Value = '3.469.28'
$ws.Range("E27").Value = '  +3.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '77.53'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.92%  '
$ws.Range("E29").Value = '  +7.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.187'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +30.05%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.07'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.99%  '
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.998'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '557.27'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.24%  '
$ws.Range("E35").Value = '  -3.14%  '
$ws.Range("E36").Value = '  -1.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.11'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +15.21%  '
$ws.Range("E38").Value = '  -9.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.91'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("B40").Value = 'WhiteBITCoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '21.78'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.62%  '
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.997'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.407'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.84%  '
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.02'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.20%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '154.38'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '182.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.56%  '
$ws.Range("B48").Value = 'ImmutableX'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.37'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.74%  '
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.19'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.94%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.31'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.757'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.60%  '
